## test/test.xlsx cleanup:
##  - rename the sheet from the LibreOffice-era "test" to Excel's default "Sheet1"
##  - re-apply the cell format on the header/data cells so the workbook carries an
##    explicit style (mirrors what Excel writes out when it touches A1/B:B6 on save)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet.
$ws.Name = "Sheet1"

# Touch the used cells (A1 and B1:B6) so Excel records an explicit,
# second cell style for them instead of implicitly inheriting style 0.
$header = $ws.Range("A1")
$body = $ws.Range("B1:B6")

$header.FormulaHidden = $true
$body.FormulaHidden = $true
